# Generate Report for Handback
#
# The handback for 9dced597-287e-4541-a4e2-02146d0d2a6d.md has now completed
# successfully for both target locales, so the generated status report is
# refreshed to reflect the new state:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - "Latest Handback DateTime" is refreshed to the new handback timestamp
#   - The stale "version mismatch" Error Detail is cleared

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$handedBackStatus = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status columns for the 9dced597 row
$wsOverview.Range("E3").Value = $handedBackStatus
$wsOverview.Range("F3").Value = $handedBackStatus

# zh-cn detail sheet, row for 9dced597
$wsZhCn.Range("C3").Value = $handedBackStatus
$wsZhCn.Range("K3").Value = "2016-09-01 09:02:43"
$wsZhCn.Range("P3").Value = ""

# de-de detail sheet, row for 9dced597
$wsDeDe.Range("C3").Value = $handedBackStatus
$wsDeDe.Range("K3").Value = "2016-09-01 09:02:50"
$wsDeDe.Range("P3").Value = ""
